$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "% of Flaky Tests" header cell to use a fresh shared-string slot
$ws.Range("D1").Value = "% of Flaky Tests"

# Relabel the test setups (order matters for shared-string insertion order)
$ws.Range("A4").Value = "Router"
$ws.Range("A2").Value = "APSNG"
$ws.Range("A3").Value = "Dave2"
$ws.Range("A5").Value = "R1"
$ws.Range("A6").Value = "R2"
$ws.Range("A7").Value = "R3"
$ws.Range("A8").Value = "R4"

# Update the active selection shown in the sheet view
$ws.Range("A9").Select()
